$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<delre>"
$ws.Range("C2").Value = 46

$ws.Range("B3").Value = "<sente>"
$ws.Range("C3").Value = 48

$ws.Range("B4").Value = "<sentence>"

$ws.Range("B5").Value = "<a>"
$ws.Range("C5").Value = 52

$ws.Range("B6").Value = "<escar>"
$ws.Range("C6").Value = 47

$ws.Range("C7").Value = 43

$ws.Range("C8").Value = 35

$ws.Range("B9").Value = "<it>"
$ws.Range("C9").Value = 14
